# Data Overview 1.xlsx edit
# - Part 1: insert "+CPI colitis C6" before the existing C7/C8 values,
#   shifting the old C7 -> C8 and old C8 -> C9.
# - Update each sheet's selection, and make "Part 4 (14 Samples)" the
#   active (tab-selected) sheet.

$wb = $excel.ActiveWorkbook

# --- Part 1 (22 Samples): shift colitis-C7/C8 labels down and insert C6 ---
$ws1 = $wb.Worksheets.Item("Part 1 (22 Samples)")
$ws1.Range("C9").Value = $ws1.Range("C8").Value2
$ws1.Range("C8").Value = $ws1.Range("C7").Value2
$ws1.Range("C7").Value = '"+CPI colitis C6"'
$ws1.Activate()
$ws1.Range("C8").Select()

# --- Part 2 (30 Samples): selection only ---
$ws2 = $wb.Worksheets.Item("Part 2 (30 Samples)")
$ws2.Activate()
$ws2.Range("B6").Select()

# --- Part 3 (31 Samples): selection only (no longer the active tab) ---
$ws3 = $wb.Worksheets.Item("Part 3 (31 Samples)")
$ws3.Activate()
$ws3.Range("B8").Select()

# --- Part 4 (14 Samples): becomes the active/tab-selected sheet ---
$ws4 = $wb.Worksheets.Item("Part 4 (14 Samples)")
$ws4.Activate()
$ws4.Range("B16").Select()
